$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Values are written with a leading apostrophe to force text storage (matching
# the original inline-string cells) and the style is reset to Normal so no
# stray quote-prefix formatting is introduced.

$ws.Range("D2").Value = "'64.801.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.68%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.157.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.94%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'571.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.13%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'150.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +4.45%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.05%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.156.11"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.96%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +4.42%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +5.35%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.08%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.506"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +7.14%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +12.52%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'38.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +8.41%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.668.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.03%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'64.852.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.72%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'7.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +7.09%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.157.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.81%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.40%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'517.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +6.57%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +6.90%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.737"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +8.90%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'15.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +6.56%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'7.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +3.97%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'85.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +4.91%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.07%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +4.40%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E29").Value = "'  +6.47%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'27.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +5.95%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.13%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +8.17%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +3.96%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'6.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +8.83%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'6.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +5.61%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'55.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.01%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'487.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +8.03%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0867"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +6.05%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +2.95%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.60%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.117.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +5.00%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'8.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +5.30%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +5.83%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.295"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +12.54%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +14.51%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'29.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +3.63%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0₃0578"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +11.81%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.02%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +3.01%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +10.25%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'118.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.29%  "
$ws.Range("E51").Style = "Normal"
